# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table (rows 16-27) lists monthly arrears periods with
# their corresponding arrears value (column F). The value that was recorded
# against period 2311 (row 16) and period 2212 (row 27) were swapped; this
# corrects them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 -> period 2311: valor mora corrected from 84000 to 90000
$ws.Range("F16").Value = 90000

# Row 27 -> period 2212: valor mora corrected from 90000 to 84000
$ws.Range("F27").Value = 84000
